# Kayıt silindi: 11117797
# The record with Kayıt No 11117797 (CİNS DEĞ. / HAVVA NİLGÜN KIYMAÇ, ÖZLEM AYDINLI)
# dated 2025-11-11 is removed from both the master "Kayitlar" sheet (row 612)
# and its per-district mirror "Merkez İlçe" sheet (row 70). Deleting the
# entire row shifts every following row up by one and shrinks each sheet's
# used range by one row, exactly matching the target diff.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(612).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(70).Delete()
